$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 2-10: full Il4/Il2rg ligand-receptor data, now with "ECs" added as a
# sending cluster (3 senders x 3 targets = 9 rows, was 2 senders x 3 targets = 6 rows).

# Row 2: ECs -> ECs
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Il4"
$ws.Range("C2").Value = "Il2rg"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.6580480000000001
$ws.Range("H2").Value = 1.974144
$ws.Range("I2").Value = 0.1830904640197835
$ws.Range("J2").Value = 0.1830904640197835
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 32.30682666666667
$ws.Range("N2").Value = 96.92048
$ws.Range("O2").Value = 0.886587237369156
$ws.Range("P2").Value = 0.8865872373691559
$ws.Range("Q2").Value = 21.25944267434667
$ws.Range("R2").Value = 191.33498406912
$ws.Range("S2").Value = 0.1623256686839367
$ws.Range("T2").Value = 0.1623256686839367

# Row 3: ECs -> FAPs
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Il4"
$ws.Range("C3").Value = "Il2rg"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 0.6580480000000001
$ws.Range("H3").Value = 1.974144
$ws.Range("I3").Value = 0.1830904640197835
$ws.Range("J3").Value = 0.1830904640197835
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 0.7369563333333332
$ws.Range("N3").Value = 2.210869
$ws.Range("O3").Value = 0.02022408719906369
$ws.Range("P3").Value = 0.02022408719906369
$ws.Range("Q3").Value = 0.4849526412373333
$ws.Range("R3").Value = 4.364573771136
$ws.Range("S3").Value = 0.003702837509653133
$ws.Range("T3").Value = 0.003702837509653133

# Row 4: ECs -> sCs
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Il4"
$ws.Range("C4").Value = "Il2rg"
$ws.Range("D4").Value = "sCs"
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 0.6580480000000001
$ws.Range("H4").Value = 1.974144
$ws.Range("I4").Value = 0.1830904640197835
$ws.Range("J4").Value = 0.1830904640197835
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 3.395752
$ws.Range("N4").Value = 10.187256
$ws.Range("O4").Value = 0.09318867543178035
$ws.Range("P4").Value = 0.09318867543178033
$ws.Range("Q4").Value = 2.234567812096
$ws.Range("R4").Value = 20.111110308864
$ws.Range("S4").Value = 0.01706195782619366
$ws.Range("T4").Value = 0.01706195782619365

# Row 5: FAPs -> ECs
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Il4"
$ws.Range("C5").Value = "Il2rg"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 1.115367
$ws.Range("H5").Value = 3.346101
$ws.Range("I5").Value = 0.3103315587652478
$ws.Range("J5").Value = 0.3103315587652479
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 32.30682666666667
$ws.Range("N5").Value = 96.92048
$ws.Range("O5").Value = 0.886587237369156
$ws.Range("P5").Value = 0.8865872373691559
$ws.Range("Q5").Value = 36.03396833872
$ws.Range("R5").Value = 324.30571504848
$ws.Range("S5").Value = 0.275135999354145
$ws.Range("T5").Value = 0.275135999354145

# Row 6: FAPs -> FAPs
$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Il4"
$ws.Range("C6").Value = "Il2rg"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 1.115367
$ws.Range("H6").Value = 3.346101
$ws.Range("I6").Value = 0.3103315587652478
$ws.Range("J6").Value = 0.3103315587652479
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 0.7369563333333332
$ws.Range("N6").Value = 2.210869
$ws.Range("O6").Value = 0.02022408719906369
$ws.Range("P6").Value = 0.02022408719906369
$ws.Range("Q6").Value = 0.8219767746409998
$ws.Range("R6").Value = 7.397790971768999
$ws.Range("S6").Value = 0.006276172505089729
$ws.Range("T6").Value = 0.00627617250508973

# Row 7: FAPs -> sCs
$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Il4"
$ws.Range("C7").Value = "Il2rg"
$ws.Range("D7").Value = "sCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 1.115367
$ws.Range("H7").Value = 3.346101
$ws.Range("I7").Value = 0.3103315587652478
$ws.Range("J7").Value = 0.3103315587652479
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 3.395752
$ws.Range("N7").Value = 10.187256
$ws.Range("O7").Value = 0.09318867543178035
$ws.Range("P7").Value = 0.09318867543178033
$ws.Range("Q7").Value = 3.787509720984
$ws.Range("R7").Value = 34.087587488856
$ws.Range("S7").Value = 0.02891938690601315
$ws.Range("T7").Value = 0.02891938690601315

# Row 8: sCs -> ECs
$ws.Range("A8").Value = "sCs"
$ws.Range("B8").Value = "Il4"
$ws.Range("C8").Value = "Il2rg"
$ws.Range("D8").Value = "ECs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 1.820699
$ws.Range("H8").Value = 5.462097
$ws.Range("I8").Value = 0.5065779772149687
$ws.Range("J8").Value = 0.5065779772149687
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 32.30682666666667
$ws.Range("N8").Value = 96.92048
$ws.Range("O8").Value = 0.886587237369156
$ws.Range("P8").Value = 0.8865872373691559
$ws.Range("Q8").Value = 58.82100700517334
$ws.Range("R8").Value = 529.38906304656
$ws.Range("S8").Value = 0.4491255693310743
$ws.Range("T8").Value = 0.4491255693310743

# Row 9: sCs -> FAPs
$ws.Range("A9").Value = "sCs"
$ws.Range("B9").Value = "Il4"
$ws.Range("C9").Value = "Il2rg"
$ws.Range("D9").Value = "FAPs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 1.820699
$ws.Range("H9").Value = 5.462097
$ws.Range("I9").Value = 0.5065779772149687
$ws.Range("J9").Value = 0.5065779772149687
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 0.7369563333333332
$ws.Range("N9").Value = 2.210869
$ws.Range("O9").Value = 0.02022408719906369
$ws.Range("P9").Value = 0.02022408719906369
$ws.Range("Q9").Value = 1.341775659143666
$ws.Range("R9").Value = 12.075980932293
$ws.Range("S9").Value = 0.01024507718432082
$ws.Range("T9").Value = 0.01024507718432082

# Row 10: sCs -> sCs
$ws.Range("A10").Value = "sCs"
$ws.Range("B10").Value = "Il4"
$ws.Range("C10").Value = "Il2rg"
$ws.Range("D10").Value = "sCs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 1.820699
$ws.Range("H10").Value = 5.462097
$ws.Range("I10").Value = 0.5065779772149687
$ws.Range("J10").Value = 0.5065779772149687
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 3.395752
$ws.Range("N10").Value = 10.187256
$ws.Range("O10").Value = 0.09318867543178035
$ws.Range("P10").Value = 0.09318867543178033
$ws.Range("Q10").Value = 6.182642270648
$ws.Range("R10").Value = 55.643780435832
$ws.Range("S10").Value = 0.04720733069957354
$ws.Range("T10").Value = 0.04720733069957353

Write-Host "done"
